$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing row 67 values (3rd data correction before appending new row)
$ws.Range("B67").Value = 125.5
$ws.Range("C67").Value = 97.09999999999999
$ws.Range("D67").Value = 125.9

# Append new row 68 with the new monthly data point
$ws.Range("A68").Value = "'01-07-2021"
$ws.Range("A67").Copy()
$ws.Range("A68").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B68").Value = 127.1
$ws.Range("C68").Value = 97.59999999999999
$ws.Range("D68").Value = 127.7
